$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 99.99992847442627
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 99.9825656414032
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 60.59531569480896
$ws.Range("D9").Value = 100
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 85.95828413963318
$ws.Range("C11").Value = 2
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 91.99870824813843
$ws.Range("D15").Value = 100
$ws.Range("D16").Value = 100
$ws.Range("D17").Value = 97.26511240005493
$ws.Range("D20").Value = 99.99951124191284
$ws.Range("D21").Value = 99.99995231628418
$ws.Range("D22").Value = 100
$ws.Range("D25").Value = 100
$ws.Range("D26").Value = 99.99978542327881
$ws.Range("D30").Value = 99.99996423721313
$ws.Range("D32").Value = 99.99543428421021
$ws.Range("D35").Value = 100
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 100
$ws.Range("D38").Value = 99.99998807907104
$ws.Range("D43").Value = 99.72549080848694
$ws.Range("D46").Value = 100
$ws.Range("D48").Value = 100
$ws.Range("D54").Value = 100
$ws.Range("D55").Value = 99.99998807907104
$ws.Range("D56").Value = 99.96091723442078
$ws.Range("D57").Value = 99.9971866607666
$ws.Range("D59").Value = 100
$ws.Range("C61").Value = 1
$ws.Range("D62").Value = 100
$ws.Range("D63").Value = 99.99991655349731
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 99.98040795326233
$ws.Range("D69").Value = 99.99990463256836
$ws.Range("D70").Value = 100
$ws.Range("D72").Value = 100
$ws.Range("D73").Value = 99.99181032180786
$ws.Range("D78").Value = 99.99961853027344
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = 84.98300909996033
$ws.Range("D80").Value = 100
$ws.Range("D81").Value = 100
$ws.Range("C82").Value = 1
$ws.Range("D82").Value = 99.99998807907104
$ws.Range("D84").Value = 99.99997615814209
$ws.Range("D87").Value = 100
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = 53.71508002281189
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = 100
$ws.Range("C90").Value = 2
$ws.Range("D90").Value = 100
$ws.Range("D91").Value = 100
$ws.Range("D92").Value = 99.99669790267944
$ws.Range("D93").Value = 99.89771842956543
$ws.Range("D94").Value = 89.76674675941467
$ws.Range("D96").Value = 99.99967813491821
$ws.Range("D98").Value = 100
$ws.Range("D100").Value = 100
$ws.Range("D101").Value = 99.98487234115601
$ws.Range("D102").Value = 100
$ws.Range("D107").Value = 99.99212026596069
$ws.Range("D108").Value = 99.99955892562866
$ws.Range("D109").Value = 100
$ws.Range("C110").Value = 2
$ws.Range("D110").Value = 100
$ws.Range("D111").Value = 100
$ws.Range("D115").Value = 100
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 99.99998807907104
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 99.99616146087646
$ws.Range("D118").Value = 99.97789263725281
$ws.Range("D122").Value = 99.99954700469971
$ws.Range("D124").Value = 100
$ws.Range("C127").Value = 2
$ws.Range("D127").Value = 100
$ws.Range("D128").Value = 99.99996423721313
$ws.Range("D130").Value = 100
$ws.Range("C133").Value = 1
$ws.Range("D133").Value = 99.99511241912842
$ws.Range("D134").Value = 100
$ws.Range("D135").Value = 99.99905824661255
$ws.Range("D136").Value = 100
$ws.Range("D140").Value = 99.98816251754761
$ws.Range("D143").Value = 99.99986886978149
$ws.Range("D144").Value = 100
$ws.Range("D153").Value = 100
$ws.Range("D155").Value = 100
$ws.Range("D156").Value = 100
$ws.Range("D157").Value = 100
$ws.Range("D160").Value = 100
$ws.Range("C161").Value = 2
$ws.Range("D161").Value = 92.17342138290405
$ws.Range("D162").Value = 100
$ws.Range("D163").Value = 100
$ws.Range("D164").Value = 95.80987691879272
$ws.Range("D166").Value = 81.40982389450073
$ws.Range("D167").Value = 99.99996423721313
$ws.Range("D168").Value = 99.99998807907104
$ws.Range("D169").Value = 99.99997615814209
$ws.Range("D170").Value = 99.99955892562866
$ws.Range("D173").Value = 100
$ws.Range("D178").Value = 99.99996423721313
$ws.Range("D179").Value = 99.9994158744812
$ws.Range("D181").Value = 100
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 97.34569787979126
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 100
$ws.Range("C186").Value = 1
$ws.Range("D186").Value = 64.0005350112915
$ws.Range("D188").Value = 100
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 100
$ws.Range("C193").Value = 2
$ws.Range("D193").Value = 100
$ws.Range("D196").Value = 99.99997615814209
$ws.Range("C199").Value = 2
$ws.Range("D199").Value = 99.99889135360718
$ws.Range("D200").Value = 99.99468326568604
$ws.Range("D202").Value = 100
$ws.Range("C204").Value = 2
$ws.Range("D204").Value = 99.99995231628418
$ws.Range("D206").Value = 100
$ws.Range("D210").Value = 99.99997615814209
$ws.Range("D211").Value = 63.68331909179688
$ws.Range("D213").Value = 99.99988079071045
$ws.Range("D214").Value = 99.90843534469604
$ws.Range("D216").Value = 99.99996423721313
$ws.Range("C217").Value = 1
$ws.Range("D217").Value = 99.99675750732422
$ws.Range("D220").Value = 100
$ws.Range("C221").Value = 1
$ws.Range("D222").Value = 100
$ws.Range("D223").Value = 99.99940395355225
$ws.Range("D224").Value = 100
$ws.Range("D225").Value = 100
$ws.Range("D226").Value = 100
$ws.Range("C227").Value = 1
$ws.Range("D227").Value = 100
$ws.Range("C229").Value = 1
$ws.Range("D229").Value = 99.9998927116394
$ws.Range("D230").Value = 100
$ws.Range("D231").Value = 100
$ws.Range("C233").Value = 1
$ws.Range("D233").Value = 100

Write-Output "Applied 152 changes"
